# Restore original simultaneous logic with 16:00 exit window
#
# The trade log currently has columns: Date, Type, Entry, Exit, PnL, Result, Close
# It needs to become:                  Date, Type, Entry, Exit, Close, PnL, Result
#
# i.e. a "Close" column (duplicating the Exit price) is reinserted right after
# "Exit", pushing "PnL" and "Result" one column to the right, and the old
# trailing "Close" column (G) goes away.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at E — shifts old E (PnL), F (Result), G (Close)
# one column to the right, so the data that used to live in G now lives in H.
$ws.Columns.Item(5).Insert()

# Copy the (now shifted) "Close" column H back into the freshly inserted
# column E, then remove the now-duplicate column H.
$ws.Range("H1:H26").Copy()
$ws.Range("E1:E26").PasteSpecial()
$ws.Columns.Item(8).Delete()
